$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: "gepürft" -> "geprüft" (typo fix), dropping the spell-check
# proofErr markers and re-splitting the sentence into three runs:
# " gep" | "rü" | "ft und dann die Typen-Abkürzung gesetzt."
# -----------------------------------------------------------------

$rng = $d.Content
$rng.Find.Execute(" gepürft und dann die Typen-Abkürzung gesetzt.", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    " geprüft und dann die Typen-Abkürzung gesetzt.", 2)

$rng2 = $d.Content
$rng2.Find.Execute(" geprüft und dann die Typen-Abkürzung gesetzt.")
$start = $rng2.Start
$end = $rng2.End

# Clear the text then retype it in three separate chunks so Word
# creates three distinct runs at those exact boundaries.
$delRange = $d.Range($start, $end)
$delRange.Text = ""

$ins1 = $d.Range($start, $start)
$ins1.InsertAfter(" gep")

$ins2 = $d.Range($start + 4, $start + 4)
$ins2.InsertAfter("rü")

$ins3 = $d.Range($start + 6, $start + 6)
$ins3.InsertAfter("ft und dann die Typen-Abkürzung gesetzt.")

# -----------------------------------------------------------------
# Change 2 & 3: move the "_GoBack" bookmark from the end of the
# document (after "...Link am Ende der Zeile verbunden.") into the
# middle of the word "Reihe" (between "Rei" and "he").
# -----------------------------------------------------------------

$findReihe = $d.Content
$findReihe.Find.Execute("Reihe der den Switchen")
$reiheStart = $findReihe.Start

$bmRange = $d.Range($reiheStart + 3, $reiheStart + 3)
$d.Bookmarks.Add("_GoBack", $bmRange)
